{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Identify the run of paragraphs to remove: the empty paragraph and the\n// page-break paragraph that precede the trailing copyright notice, plus the\n// copyright paragraph itself.\nconst copyrightText =\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\";\n\nlet copyrightIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === copyrightText) {\n    copyrightIndex = i;\n    break;\n  }\n}\n\nif (copyrightIndex === -1) {\n  throw new Error(\"Could not locate the copyright paragraph to remove.\");\n}\n\n// The two empty paragraphs immediately before the copyright paragraph\n// (a plain blank paragraph followed by a page-break paragraph) are removed\n// along with the copyright paragraph itself.\nconst toDelete = [copyrightIndex, copyrightIndex - 1, copyrightIndex - 2];\ntoDelete.sort((a, b) => b - a);\nfor (const idx of toDelete) {\n  paragraphs.items[idx].delete();\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the trailing copyright/footer paragraph.\n$searchRange = $d.Content\n$searchRange.Find.ClearFormatting()\n$found = $searchRange.Find.Execute(\"Powered by Jekyll and Github pages\")\nif (-not $found) {\n    throw \"Could not locate the copyright paragraph to remove.\"\n}\n$searchRange.Expand(4) | Out-Null  # wdParagraph\n\n# Resolve the 1-based Paragraphs index of that paragraph.\n$copyrightIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Start -eq $searchRange.Start) {\n        $copyrightIndex = $i\n        break\n    }\n}\nif ($copyrightIndex -eq -1) {\n    throw \"Could not resolve paragraph index for the copyright paragraph.\"\n}\n\n# Remove the copyright paragraph plus the blank paragraph and the\n# page-break paragraph that immediately precede it, highest index first so\n# earlier indices stay valid.\n$d.Paragraphs.Item($copyrightIndex).Range.Delete()\n$d.Paragraphs.Item($copyrightIndex - 1).Range.Delete()\n$d.Paragraphs.Item($copyrightIndex - 2).Range.Delete()\n"}
